# Update progress values on the "Fortschritt" sheet and move the
# active selection, matching the authored workbook edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fortschritt")

# Mussziele block
$ws.Range("A11").Value = 1

# Optional block
$ws.Range("A15").Value = 1
$ws.Range("A16").Value = 0.7
$ws.Range("A19").Value = 1

# Leave the selection where the author left it when saving.
$ws.Range("A21").Select()
